# Auto-generated edit script: updates cryptos price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "21.095.30"
$ws.Range("E2").Value = "  -4.28%  "

# Row 3
Set-TextValue "D3" "1.508.90"
$ws.Range("E3").Value = "  -2.88%  "

# Row 4
Set-TextValue "D4" "1.007"
$ws.Range("E4").Value = "  +0.56%  "

# Row 5
Set-TextValue "D5" "1.007"
$ws.Range("E5").Value = "  +0.56%  "

# Row 6
Set-TextValue "D6" "284.81"
$ws.Range("E6").Value = "  -1.81%  "

# Row 7
Set-TextValue "D7" "0.3866"
$ws.Range("E7").Value = "  -1.80%  "

# Row 8
Set-TextValue "D8" "0.3131"
$ws.Range("E8").Value = "  -2.74%  "

# Row 9
Set-TextValue "D9" "42.23"
$ws.Range("E9").Value = "  -3.19%  "

# Row 10
Set-TextValue "D10" "0.06996"
$ws.Range("E10").Value = "  -3.23%  "

# Row 11
Set-TextValue "D11" "1.041"
$ws.Range("E11").Value = "  -2.89%  "

# Row 12
Set-TextValue "D12" "1.007"
$ws.Range("E12").Value = "  +0.63%  "

# Row 13
Set-TextValue "D13" "5.621"
$ws.Range("E13").Value = "  -0.64%  "

# Row 14
Set-TextValue "D14" "17.86"
$ws.Range("E14").Value = "  -4.53%  "

# Row 15
Set-TextValue "D15" "1.517.31"
$ws.Range("E15").Value = "  -2.23%  "

# Row 16
Set-TextValue "D16" "6.376"
$ws.Range("E16").Value = "  -3.57%  "

# Row 17
Set-TextValue "D17" "0.00001070"
$ws.Range("E17").Value = "  -5.28%  "

# Row 18
Set-TextValue "D18" "0.06576"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19
Set-TextValue "D19" "82.11"
$ws.Range("E19").Value = "  -1.45%  "

# Row 20
$ws.Range("E20").Value = "  +0.63%  "

# Row 21
Set-TextValue "D21" "6.004"
$ws.Range("E21").Value = "  -4.15%  "

# Row 22
Set-TextValue "D22" "15.16"
$ws.Range("E22").Value = "  -1.92%  "

# Row 23
Set-TextValue "D23" "10.84"
$ws.Range("E23").Value = "  -3.86%  "

# Row 24
Set-TextValue "D24" "2.351"
$ws.Range("E24").Value = "  -0.41%  "

# Row 25
Set-TextValue "D25" "21.102.93"
$ws.Range("E25").Value = "  -4.30%  "

# Row 26
Set-TextValue "D26" "2.346"
$ws.Range("E26").Value = "  -2.73%  "

# Row 27
Set-TextValue "D27" "147.87"
$ws.Range("E27").Value = "  -0.73%  "

# Row 28
Set-TextValue "D28" "18.04"
$ws.Range("E28").Value = "  -2.59%  "

# Row 29
Set-TextValue "D29" "4.795"
$ws.Range("E29").Value = "  -1.68%  "

# Row 30
Set-TextValue "D30" "1.680.11"
$ws.Range("E30").Value = "  -2.65%  "

# Row 31
Set-TextValue "D31" "114.77"
$ws.Range("E31").Value = "  -3.02%  "

# Row 32
Set-TextValue "D32" "5.938"
$ws.Range("E32").Value = "  +2.28%  "

# Row 33
Set-TextValue "D33" "0.9582"
$ws.Range("E33").Value = "  -2.10%  "

# Row 34
Set-TextValue "D34" "0.07954"
$ws.Range("E34").Value = "  -4.41%  "

# Row 35
Set-TextValue "D35" "8.411"
$ws.Range("E35").Value = "  -7.12%  "

# Row 36
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D36" "5.081"
$ws.Range("E36").Value = "  -0.46%  "

# Row 37
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D37" "11.29"
$ws.Range("E37").Value = "  +5.95%  "

# Row 38
$ws.Range("B38").Value = "WEMIXTOKEN"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D38" "1.467"
$ws.Range("E38").Value = "  -8.34%  "

# Row 39
Set-TextValue "D39" "0.05800"
$ws.Range("E39").Value = "  -3.37%  "

# Row 40
Set-TextValue "D40" "0.02129"
$ws.Range("E40").Value = "  -5.76%  "

# Row 41
Set-TextValue "D41" "1.006"
$ws.Range("E41").Value = "  +0.57%  "

# Row 42
Set-TextValue "D42" "1.158"
$ws.Range("E42").Value = "  -4.08%  "

# Row 43
Set-TextValue "D43" "0.1974"
$ws.Range("E43").Value = "  -2.85%  "

# Row 44
Set-TextValue "D44" "0.5618"
$ws.Range("E44").Value = "  -3.02%  "

# Row 45
Set-TextValue "D45" "12.92"
$ws.Range("E45").Value = "  +0.39%  "

# Row 46
Set-TextValue "D46" "3.681"
$ws.Range("E46").Value = "  -1.58%  "

# Row 47
Set-TextValue "D47" "0.5431"
$ws.Range("E47").Value = "  -2.31%  "

# Row 48
$ws.Range("E48").Value = "  +0.61%  "

# Row 49
Set-TextValue "D49" "1.849"
$ws.Range("E49").Value = "  -2.24%  "

# Row 50
Set-TextValue "D50" "113.72"
$ws.Range("E50").Value = "  -3.31%  "

# Row 51
Set-TextValue "D51" "0.06561"
$ws.Range("E51").Value = "  -3.75%  "
